$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 5422.143
$ws.Range("I88").Value = 553.25
$ws.Range("J88").Value = 8418.385
$ws.Range("K88").Value = 553.25
$ws.Range("L88").Value = 8418.385
$ws.Range("M88").Value = -147.25
$ws.Range("N88").Value = -9230.385

$ws.Range("H91").Value = 5422.143
$ws.Range("I91").Value = 553.25
$ws.Range("J91").Value = 8418.385
$ws.Range("K91").Value = 553.25
$ws.Range("L91").Value = 8418.385
$ws.Range("M91").Value = 850.75
$ws.Range("N91").Value = -11226.385

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2517.52
$ws.Range("I32").Value = 2134.4
$ws.Range("J32").Value = 5965.6
$ws.Range("K32").Value = 2134.4
$ws.Range("L32").Value = 5965.6
$ws.Range("M32").Value = -1847.4
$ws.Range("N32").Value = -6539.6

$ws.Range("M54").ClearContents()
$ws.Range("H54").Value = 17514
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 17514
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 17514
$ws.Range("N54").Value = -19052

$ws.Range("H61").Value = 271866.84
$ws.Range("I61").Value = 1601.0667
$ws.Range("J61").Value = 1430148.8
$ws.Range("K61").Value = 1601.0667
$ws.Range("L61").Value = 1430148.8
$ws.Range("M61").Value = -1389.0667
$ws.Range("N61").Value = -1430572.8

$ws.Range("H124").Value = 23833.334
$ws.Range("J124").Value = 23833.334
$ws.Range("L124").Value = 23833.334
$ws.Range("N124").Value = -33653.334

$ws.Range("H132").Value = 6078.7046
$ws.Range("I132").Value = 5328.64
$ws.Range("J132").Value = 7065.6313
$ws.Range("K132").Value = 15985.92
$ws.Range("L132").Value = 21196.8939
$ws.Range("M132").Value = -13455.92
$ws.Range("N132").Value = -26256.8939

$ws.Range("H136").Value = 271866.84
$ws.Range("I136").Value = 1601.0667
$ws.Range("J136").Value = 1430148.8
$ws.Range("K136").Value = 4803.2001
$ws.Range("L136").Value = 4290446.4
$ws.Range("M136").Value = -2253.2001
$ws.Range("N136").Value = -4295546.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H98").Value = 48000
$ws.Range("J98").Value = 48000
$ws.Range("L98").Value = 48000
$ws.Range("N98").Value = -53990

$ws.Range("H107").Value = 3102.077
$ws.Range("I107").Value = 2683.8572
$ws.Range("J107").Value = 3590
$ws.Range("K107").Value = 2683.8572
$ws.Range("L107").Value = 3590
$ws.Range("M107").Value = -763.8571999999999
$ws.Range("N107").Value = -7430

$ws.Range("H122").Value = 30780
$ws.Range("J122").Value = 30780
$ws.Range("L122").Value = 30780
$ws.Range("N122").Value = -40580

$ws.Range("H125").Value = 37868
$ws.Range("I125").Value = 40000
$ws.Range("J125").Value = 37335
$ws.Range("K125").Value = 40000
$ws.Range("L125").Value = 37335
$ws.Range("M125").Value = -35080
$ws.Range("N125").Value = -47175

$ws.Range("H130").Value = 37526.668
$ws.Range("J130").Value = 37526.668
$ws.Range("L130").Value = 37526.668
$ws.Range("N130").Value = -47566.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 11378189
$ws.Range("I31").Value = 20834470
$ws.Range("J31").Value = 30650.3
$ws.Range("K31").Value = 20834470
$ws.Range("L31").Value = 30650.3
$ws.Range("M31").Value = -20834175
$ws.Range("N31").Value = -31240.3

$ws.Range("H34").Value = 11378189
$ws.Range("I34").Value = 20834470
$ws.Range("J34").Value = 30650.3
$ws.Range("K34").Value = 20834470
$ws.Range("L34").Value = 30650.3
$ws.Range("M34").Value = -20834268
$ws.Range("N34").Value = -31054.3

$ws.Range("H107").Value = 53624.316
$ws.Range("I107").Value = 67670.664
$ws.Range("J107").Value = 950.5
$ws.Range("K107").Value = 67670.664
$ws.Range("L107").Value = 950.5
$ws.Range("M107").Value = -65750.664
$ws.Range("N107").Value = -4790.5

$ws.Range("H124").Value = 21471
$ws.Range("J124").Value = 21471
$ws.Range("L124").Value = 21471
$ws.Range("N124").Value = -26381

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 135.04347
$ws.Range("I33").Value = 34
$ws.Range("J33").Value = 498.8
$ws.Range("K33").Value = 204
$ws.Range("L33").Value = 2992.8
$ws.Range("M33").Value = 79
$ws.Range("N33").Value = -3558.8

$ws.Range("H76").Value = 2092.2856
$ws.Range("I76").Value = 882
$ws.Range("K76").Value = 2646
$ws.Range("M76").Value = -2263

$ws.Range("H79").Value = 2092.2856
$ws.Range("I79").Value = 882
$ws.Range("K79").Value = 2646
$ws.Range("M79").Value = -1320

$ws.Range("H88").Value = 4909.091
$ws.Range("J88").Value = 4909.091
$ws.Range("L88").Value = 14727.273
$ws.Range("N88").Value = -15583.273

$ws.Range("H91").Value = 4909.091
$ws.Range("J91").Value = 4909.091
$ws.Range("L91").Value = 14727.273
$ws.Range("N91").Value = -17691.273

$ws.Range("H136").Value = 1787.5217
$ws.Range("I136").Value = 975.8125
$ws.Range("J136").Value = 3642.8572
$ws.Range("K136").Value = 2927.4375
$ws.Range("L136").Value = 10928.5716
$ws.Range("M136").Value = 2172.5625
$ws.Range("N136").Value = -21128.5716

$ws.Range("N138").ClearContents()
$ws.Range("H138").Value = 3260
$ws.Range("I138").Value = 3260
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 9780
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = -4640

$ws.Range("H139").Value = 689228.0600000001
$ws.Range("I139").Value = 734976.6
$ws.Range("J139").Value = 3000
$ws.Range("K139").Value = 2204929.8
$ws.Range("L139").Value = 9000
$ws.Range("M139").Value = -2199789.8
$ws.Range("N139").Value = -19280

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 5000
$ws.Range("J4").Value = 5000
$ws.Range("L4").Value = 5000
$ws.Range("N4").Value = -5224

$ws.Range("H107").Value = 243.625
$ws.Range("I107").Value = 266.5
$ws.Range("J107").Value = 175
$ws.Range("K107").Value = 266.5
$ws.Range("L107").Value = 175
$ws.Range("M107").Value = 1653.5
$ws.Range("N107").Value = -4015
